# UnitBalanceSpreadsheet.xlsx edits
# - Add "Testing results" column header (K1) and a "Ghoul" unit row (row 5)
# - Re-balance the Archer row (row 4)
# - Re-point/normalize the J-column DPS formula to use the /100 scaling factor
# - Cosmetic: row 1 height, column K width, and active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column K header: "Testing results" ---
$ws.Range("K1").Value = "Testing results"

# --- Row 1 height: 60 -> 45 ---
$ws.Rows.Item(1).RowHeight = 45

# --- Column K (11) width -> ~27.43 characters (custom width) ---
$ws.Columns.Item(11).ColumnWidth = 26.67

# --- Row 4 (Archer) data tweaks ---
$ws.Range("B4").Value = 0.6
$ws.Range("G4").Value = 75

# --- Row 5: new "Ghoul" unit ---
$ws.Range("A5").Value = "Ghoul"
$ws.Range("B5").Value = 1.5
$ws.Range("C5").Value = 10
$ws.Range("G5").Value = 80

# --- J3:J33 DPS formula now consistently uses the /100 scaling factor ---
$ws.Range("J3:J33").Formula = "=D3*F3/10 +G3*I3/100"

# --- Update the saved active cell/selection to G5 ---
[void]$ws.Range("G5").Select()
